$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete entire row 81 (the "キムネコウヨウジャク。" post) and shift all rows below up by one.
$ws.Rows.Item(81).Delete()
